# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.127
$ws.Range("E2").Value = -0.09480000000000001
$ws.Range("G2").Value = 0.05661306703513689
$ws.Range("H2").Value = 0.05661306703513689
$ws.Range("I2").Value = 0.05024039195002986
$ws.Range("J2").Value = 0.0405461819614476
$ws.Range("K2").Value = 735.7
$ws.Range("L2").Value = 0.03206824254523422
$ws.Range("M2").Value = 186.2
$ws.Range("N2").Value = 0.04246294184720639
$ws.Range("O2").Value = 0.253092293054234
$ws.Range("P2").Value = 186.2
$ws.Range("Q2").Value = 0.04246294184720639
$ws.Range("R2").Value = 0.253092293054234
$ws.Range("U2").Value = 1896.6
$ws.Range("V2").Value = 0.4325199543899658
$ws.Range("W2").Value = 0.06127310130008579
$ws.Range("X2").Value = 0.1648961547883717
$ws.Range("Y2").Value = -0.103623053488286
$ws.Range("Z2").Value = 1.537358940681373
$ws.Range("AA2").Value = 0.06233403534892529
$ws.Range("AB2").Value = 0.07330197552344268
$ws.Range("AC2").Value = -0.0109679401745174
$ws.Range("AD2").Value = 7604.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 7604.4
$ws.Range("AG2").Value = 5707.799999999999
$ws.Range("AH2").Value = 0.634260263232522
$ws.Range("AI2").Value = 0.3536946683472946
$ws.Range("AJ2").Value = 0.5655318642993025
$ws.Range("AK2").Value = 0.2911652629914351
$ws.Range("AL2").Value = 402.6
$ws.Range("AM2").Value = 402.6
$ws.Range("AN2").Value = 6.178420539486512
$ws.Range("AO2").Value = 2.862891207153502
$ws.Range("AP2").Value = 4.637471563210919
$ws.Range("AQ2").Value = 2.862891207153502

# Row 3
$ws.Range("D3").Value = 0.127
$ws.Range("E3").Value = -0.09480000000000001
$ws.Range("G3").Value = 0.05661306703513689
$ws.Range("H3").Value = 0.05661306703513689
$ws.Range("I3").Value = 0.05024039195002986
$ws.Range("J3").Value = 0.0405461819614476
$ws.Range("K3").Value = 735.7
$ws.Range("L3").Value = 0.03206824254523422
$ws.Range("M3").Value = 186.2
$ws.Range("N3").Value = 0.04246294184720639
$ws.Range("O3").Value = 0.253092293054234
$ws.Range("P3").Value = 186.2
$ws.Range("Q3").Value = 0.04246294184720639
$ws.Range("R3").Value = 0.253092293054234
$ws.Range("U3").Value = 1896.6
$ws.Range("V3").Value = 0.4325199543899658
$ws.Range("W3").Value = 0.06127310130008579
$ws.Range("X3").Value = 0.1648961547883717
$ws.Range("Y3").Value = -0.103623053488286
$ws.Range("Z3").Value = 1.537358940681373
$ws.Range("AA3").Value = 0.06233403534892529
$ws.Range("AB3").Value = 0.07330197552344268
$ws.Range("AC3").Value = -0.0109679401745174
$ws.Range("AD3").Value = 7604.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 7604.4
$ws.Range("AG3").Value = 5707.799999999999
$ws.Range("AH3").Value = 0.634260263232522
$ws.Range("AI3").Value = 0.3536946683472946
$ws.Range("AJ3").Value = 0.5655318642993025
$ws.Range("AK3").Value = 0.2911652629914351
$ws.Range("AL3").Value = 402.6
$ws.Range("AM3").Value = 402.6
$ws.Range("AN3").Value = 6.178420539486512
$ws.Range("AO3").Value = 2.862891207153502
$ws.Range("AP3").Value = 4.637471563210919
$ws.Range("AQ3").Value = 2.862891207153502
